# Apply the "Row Reader" data edits to the "Week" sheet and update the
# active sheet / selection to reflect where the author left off.

$wb = $excel.ActiveWorkbook

$week = $wb.Worksheets.Item("Week")

# Row 2 (John): toggle C2/E2/G2 to "X" and D2 back to "I"
$week.Range("C2").Value = "X"
$week.Range("D2").Value = "I"
$week.Range("E2").Value = "X"
$week.Range("G2").Value = "X"

# Row 3 (Jeffery): toggle D3/E3 to "X"
$week.Range("D3").Value = "X"
$week.Range("E3").Value = "X"

# Row 5 (Julia): reset the whole row back to "I"
$week.Range("B5").Value = "I"
$week.Range("C5").Value = "I"
$week.Range("D5").Value = "I"
$week.Range("E5").Value = "I"
$week.Range("F5").Value = "I"
$week.Range("G5").Value = "I"

# Leave the workbook with "Week" as the active sheet and cell E3 selected,
# matching where the author was last working.
$week.Activate()
$week.Range("E3").Select()
